$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A2")

# Update A2 value to the new email text (replaces the old numeric value)
$r.Value = "venkatshamuthu@gmail.com"

# Re-point the existing hyperlink at the same address so it loses its
# stale "display" override (which no longer matches the cell's text)
# while keeping the same target / relationship and cell style.
$r.Hyperlinks.Delete()
$h = $r.Hyperlinks.Item(1)
$h.Address = "mailto:venkatsoumuthu@gmail.com"

# Move the active selection to A2
$r.Select()
